# Updates the cryptocurrency price (column D) and 1h-volume-change (column E)
# cells, matching the scraper refresh captured in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new text value.
$updates = [ordered]@{
    "D2" = "64.349.82"
    "E2" = "  +1.42%  "
    "D3" = "3.174.57"
    "E3" = "  +2.49%  "
    "E4" = "  +0.09%  "
    "D5" = "593.66"
    "E5" = "  +1.93%  "
    "D6" = "148.45"
    "E6" = "  +2.24%  "
    "E7" = "  +0.02%  "
    "D8" = "3.165.61"
    "E8" = "  +2.52%  "
    "E9" = "  +1.29%  "
    "D10" = "0.162"
    "E10" = "  +1.43%  "
    "E11" = "  +5.58%  "
    "D12" = "0.462"
    "E12" = "  +1.42%  "
    "E13" = "  +1.39%  "
    "D14" = "37.69"
    "E14" = "  +1.72%  "
    "D15" = "3.699.36"
    "E15" = "  +2.48%  "
    "E16" = "  +0.19%  "
    "D17" = "7.32"
    "E17" = "  +3.86%  "
    "D18" = "64.183.26"
    "E18" = "  +1.40%  "
    "D19" = "3.162.78"
    "E19" = "  +2.17%  "
    "D20" = "471.15"
    "E20" = "  +2.39%  "
    "D21" = "14.57"
    "E21" = "  +2.79%  "
    "E22" = "  +2.49%  "
    "D23" = "7.69"
    "E23" = "  +4.00%  "
    "D24" = "2.45"
    "E24" = "  +15.00%  "
    "D25" = "13.25"
    "E25" = "  +2.91%  "
    "E26" = "  +0.60%  "
    "D27" = "10.13"
    "E27" = "  +12.39%  "
    "E28" = "  +0.04%  "
    "E29" = "  +2.85%  "
    "E30" = "  +2.82%  "
    "E31" = "  +0.14%  "
    "E32" = "  +5.55%  "
    "D33" = "28.41"
    "E33" = "  +6.83%  "
    "E34" = "  +5.79%  "
    "E35" = "  +2.19%  "
    "E36" = "  +3.41%  "
    "E37" = "  +4.13%  "
    "E38" = "  +1.38%  "
    "D39" = "3.33"
    "E39" = "  -2.00%  "
    "D40" = "469.57"
    "E40" = "  +8.74%  "
    "D41" = "51.47"
    "E41" = "  +2.46%  "
    "E42" = "  +7.70%  "
    "D43" = "0.297"
    "E43" = "  +10.50%  "
    "E44" = "  +2.96%  "
    "D45" = "2.913.30"
    "E45" = "  +1.28%  "
    "D46" = "39.85"
    "E46" = "  +11.59%  "
    "E47" = "  +0.43%  "
    "D48" = "133.22"
    "E48" = "  +6.88%  "
    "E49" = "  +0.03%  "
    "D50" = "2.26"
    "E50" = "  +5.83%  "
    "E51" = "  +1.53%  "
}

foreach ($addr in $updates.Keys) {
    $newValue = $updates[$addr]
    $cell = $ws.Range($addr)

    # These columns hold plain text in the workbook (e.g. "3.174.57" uses
    # dots as thousands separators, which Excel would otherwise happily
    # misparse as a number and round). Force the cell to Text first so the
    # literal string is preserved, then drop back to the Normal style so no
    # stray formatting is left behind on cells that used the default style.
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.Style = "Normal"
}
